$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schema")
$lo = $ws.ListObjects.Item(1)

# Add a new (blank) column to the table - it gets appended as the last
# column (F), expanding the table range from A1:E58 to A1:F58.
$newCol = $lo.ListColumns.Add(5)

# Shift the existing "Answers options" column (E) into the new slot (F),
# so column E becomes free for the new "Variable Type" column.
$ws.Range("E1:E58").Copy($ws.Range("F1:F58"))

# Write the new column header.
$ws.Range("E1").Value = "Variable Type"

# Populate "Variable Type" for each data row based on "Question type" (col B):
# Likert-scale questions are Ordinal, everything else is Nominal.
# ("Ordinal" is written before "Nominal" the first time each is used, to
#  match the shared-string intern order of the source edit.)
for ($r = 2; $r -le 58; $r++) {
    $qType = $ws.Cells.Item($r, 2).Value()
    if ($qType -eq "Likert") {
        $ws.Cells.Item($r, 5).Value = "Ordinal"
    }
}
for ($r = 2; $r -le 58; $r++) {
    $qType = $ws.Cells.Item($r, 2).Value()
    if ($qType -ne "Likert") {
        $ws.Cells.Item($r, 5).Value = "Nominal"
    }
}
